# Update "want to go" attendance counts (column F) that were refreshed
# by the gh-pages data generation job.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 408
$wsExpo.Range("F3").Value = 2363

# Sheet "演出" (performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 113

# Sheet "全部类型" (all types - combined listing)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 408
$wsAll.Range("F3").Value = 113
$wsAll.Range("F7").Value = 2363
